$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3894.2856
$ws.Range("I76").Value = 3752
$ws.Range("J76").Value = 4250
$ws.Range("K76").Value = 3752
$ws.Range("L76").Value = 4250
$ws.Range("M76").Value = -3437
$ws.Range("N76").Value = -4880
$ws.Range("H79").Value = 3894.2856
$ws.Range("I79").Value = 3752
$ws.Range("J79").Value = 4250
$ws.Range("K79").Value = 3752
$ws.Range("L79").Value = 4250
$ws.Range("M79").Value = -2660
$ws.Range("N79").Value = -6434
$ws.Range("H98").Value = 21746.592
$ws.Range("I98").Value = 23491
$ws.Range("J98").Value = 2122
$ws.Range("K98").Value = 23491
$ws.Range("L98").Value = 2122
$ws.Range("M98").Value = -21993
$ws.Range("N98").Value = -5118
$ws.Range("H122").Value = 21746.592
$ws.Range("I122").Value = 23491
$ws.Range("J122").Value = 2122
$ws.Range("K122").Value = 70473
$ws.Range("L122").Value = 6366
$ws.Range("M122").Value = -68023
$ws.Range("N122").Value = -11266
$ws.Range("H125").Value = 4173.3335
$ws.Range("I125").Value = 932
$ws.Range("J125").Value = 4821.6
$ws.Range("K125").Value = 8388
$ws.Range("L125").Value = 43394.4
$ws.Range("M125").Value = -5928
$ws.Range("N125").Value = -48314.4
$ws.Range("H129").Value = 1117.8235
$ws.Range("J129").Value = 1364.5385
$ws.Range("L129").Value = 4093.6155
$ws.Range("N129").Value = -14093.6155
$ws.Range("H134").Value = 29800
$ws.Range("J134").Value = 29800
$ws.Range("L134").Value = 29800
$ws.Range("N134").Value = -39940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1135.0667
$ws.Range("I2").Value = 1014.8182
$ws.Range("K2").Value = 1014.8182
$ws.Range("M2").Value = -901.8182
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("M12").Value = -127
$ws.Range("H14").Value = 28282.6
$ws.Range("I14").Value = 199.5
$ws.Range("J14").Value = 47004.668
$ws.Range("K14").Value = 199.5
$ws.Range("L14").Value = 47004.668
$ws.Range("M14").Value = -24.5
$ws.Range("N14").Value = -47354.668
$ws.Range("H32").Value = 6423.989
$ws.Range("I32").Value = 2835.9875
$ws.Range("J32").Value = 32518.545
$ws.Range("K32").Value = 2835.9875
$ws.Range("L32").Value = 32518.545
$ws.Range("M32").Value = -2548.9875
$ws.Range("N32").Value = -33092.545
$ws.Range("H116").Value = 1135.0667
$ws.Range("I116").Value = 1014.8182
$ws.Range("K116").Value = 1014.8182
$ws.Range("M116").Value = 1279.1818
$ws.Range("H130").Value = 28571.6
$ws.Range("J130").Value = 28571.6
$ws.Range("L130").Value = 28571.6
$ws.Range("N130").Value = -38611.6
$ws.Range("H131").Value = 77575
$ws.Range("J131").Value = 77575
$ws.Range("L131").Value = 77575
$ws.Range("N131").Value = -87655
$ws.Range("H132").Value = 1984.15
$ws.Range("I132").Value = 1452.6471
$ws.Range("J132").Value = 4996
$ws.Range("K132").Value = 4357.9413
$ws.Range("L132").Value = 14988
$ws.Range("M132").Value = -1827.9413
$ws.Range("N132").Value = -20048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1135.0667
$ws.Range("I3").Value = 1014.8182
$ws.Range("K3").Value = 1014.8182
$ws.Range("M3").Value = -900.8182
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3340
$ws.Range("H19").Value = 26669.666
$ws.Range("I19").Value = 39999
$ws.Range("J19").Value = 20005
$ws.Range("K19").Value = 39999
$ws.Range("L19").Value = 20005
$ws.Range("M19").Value = -39826
$ws.Range("N19").Value = -20351
$ws.Range("H112").Value = 31823
$ws.Range("J112").Value = 31823
$ws.Range("L112").Value = 31823
$ws.Range("N112").Value = -34777
$ws.Range("H134").Value = 960.96
$ws.Range("I134").Value = 896.6957
$ws.Range("K134").Value = 2690.0871
$ws.Range("M134").Value = -155.0870999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1596.2
$ws.Range("I94").Value = 1937.4286
$ws.Range("J94").Value = 1463.5
$ws.Range("K94").Value = 1937.4286
$ws.Range("L94").Value = 1463.5
$ws.Range("M94").Value = -1486.4286
$ws.Range("N94").Value = -2365.5
$ws.Range("H132").Value = 1978.6562
$ws.Range("I132").Value = 1735.138
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5205.414
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -2675.414
$ws.Range("N132").Value = -18057.9995
$ws.Range("H134").Value = 1610.3667
$ws.Range("I134").Value = 1438.625
$ws.Range("K134").Value = 4315.875
$ws.Range("M134").Value = -1780.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 784.5161000000001
$ws.Range("I131").Value = 357.27274
$ws.Range("J131").Value = 1019.5
$ws.Range("K131").Value = 1071.81822
$ws.Range("L131").Value = 3058.5
$ws.Range("M131").Value = 3968.18178
$ws.Range("N131").Value = -13138.5
$ws.Range("H139").Value = 3645.8333
$ws.Range("I139").Value = 1301.4166
$ws.Range("J139").Value = 5990.25
$ws.Range("K139").Value = 3904.2498
$ws.Range("L139").Value = 17970.75
$ws.Range("M139").Value = 1235.7502
$ws.Range("N139").Value = -28250.75
$ws.Range("H140").Value = 1140.4138
$ws.Range("I140").Value = 893
$ws.Range("J140").Value = 2088.8333
$ws.Range("K140").Value = 2679
$ws.Range("L140").Value = 6266.499899999999
$ws.Range("M140").Value = 2501
$ws.Range("N140").Value = -16626.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -11058
$ws.Range("H116").Value = 39312.5
$ws.Range("J116").Value = 39312.5
$ws.Range("L116").Value = 39312.5
$ws.Range("N116").Value = -48490.5
$ws.Range("H122").Value = 2655.3076
$ws.Range("I122").Value = 2739.4
$ws.Range("J122").Value = 2602.75
$ws.Range("K122").Value = 8218.200000000001
$ws.Range("L122").Value = 7808.25
$ws.Range("M122").Value = -5768.200000000001
$ws.Range("N122").Value = -12708.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2545.35
$ws.Range("I40").Value = 2505.6316
$ws.Range("K40").Value = 2505.6316
$ws.Range("M40").Value = -2369.6316
$ws.Range("H132").Value = 2081.5483
$ws.Range("I132").Value = 1057.421
$ws.Range("K132").Value = 3172.263
$ws.Range("M132").Value = -642.2629999999999
$ws.Range("H136").Value = 2198.08
$ws.Range("I136").Value = 1072.7142
$ws.Range("J136").Value = 3630.3635
$ws.Range("K136").Value = 3218.1426
$ws.Range("L136").Value = 10891.0905
$ws.Range("M136").Value = -668.1425999999997
$ws.Range("N136").Value = -15991.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7359.8
$ws.Range("J41").Value = 7359.8
$ws.Range("L41").Value = 7359.8
$ws.Range("N41").Value = -8139.8
$ws.Range("H45").Value = 5457.1113
$ws.Range("J45").Value = 5457.1113
$ws.Range("L45").Value = 5457.1113
$ws.Range("N45").Value = -6439.1113

Write-Output "Edit applied successfully"
